# Update "想去人数" (want-to-go count) figures in column F on both the
# "展览" and "全部类型" sheets, as published by the latest site regeneration.

$wb = $excel.ActiveWorkbook

$updates = @{
    "F2"  = 7432
    "F3"  = 7422
    "F4"  = 99
    "F10" = 133
    "F14" = 554
    "F17" = 4
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
